# Rename Sheet2 to "performance" and populate it with end-to-end test data.
$wb = $excel.ActiveWorkbook

# Sheet2 -> performance
$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Name = "performance"

# Populate cells column-by-column (A, B, D, C) to mirror the authoring order
$sheet2.Range("A1").Value = "job title"
$sheet2.Range("A2").Value = "Account Assistant"

$sheet2.Range("B1").Value = "kpiforjobtitle"
$sheet2.Range("B2").Value = "passedtest"

$sheet2.Range("D1").Value = "maxrating"
$sheet2.Range("D2").Value = 100

$sheet2.Range("C1").Value = "minrating"
$sheet2.Range("C2").Value = 10

# Column widths to match target
$sheet2.Columns.Item(1).ColumnWidth = 18.28515625
$sheet2.Columns.Item(2).ColumnWidth = 18.5703125
$sheet2.Columns.Item(3).ColumnWidth = 14.140625
$sheet2.Columns.Item(4).ColumnWidth = 15

# Update selection on sheet1 from B11 to B7 and deselect its tab
$sheet1 = $wb.Worksheets.Item("editorganisation")
[void]$sheet1.Range("B7").Select()

# Make performance sheet the active sheet/tab
[void]$sheet2.Select()
[void]$sheet2.Range("C2").Select()
